$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 8 for "Electricity" (added first so the new shared strings
# land in the same order as the target workbook)
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Electricity"
$ws.Range("C8").Value = "Electricity (for heat and hydrogen/synthetic fuel production)"

# Update existing row 7, column C (description) for "Other"
$ws.Range("C7").Value = "Other (incl. nuclear)"

# Update selection to match the after-state (C7 selected)
$ws.Range("C7").Select()
